$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A20").Value = 19
$ws.Range("B20").Value = 25
$ws.Range("C20").Value = 5
$ws.Range("D20").Value = 7
$ws.Range("E20").Value = 53
$ws.Range("F20").Value = 37
$ws.Range("G20").Value = 90
